$wb = $excel.ActiveWorkbook

# --- Sheet references (tab order: survey, settings, model, choices,
#     properties, queries, initial, prompt_types) ---
$survey  = $wb.Worksheets.Item("survey")
$queries = $wb.Worksheets.Item("queries")

# --- Text content changes: "individual(s)" -> "member(s)" wording ---

# survey sheet: update the household-members prompt text
$survey.Range("E16").Value = "Make a list of all members who normally live in this household"

# queries sheet: rename custom_individuals -> custom_members
$queries.Range("C2").Value = "custom_members"
$queries.Range("D2").Value = "custom_members"
$queries.Range("C3").Value = "custom_members"
$queries.Range("D3").Value = "custom_members"

# --- View / selection state: active tab moves from "survey" to "queries" ---

# Make "queries" the active/selected sheet (sets workbookView activeTab
# and the sheet's tabSelected attribute; clears it from "survey").
$queries.Activate()
$queries.Range("D8").Select() | Out-Null

# Update the lingering selection on the survey sheet too.
$survey.Range("E16").Select() | Out-Null
$survey.Activate()

# Re-activate queries so it ends up as the final active sheet.
$queries.Activate()
